$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1725.8572
$ws.Range("I55").Value = 5149.5
$ws.Range("J55").Value = 356.4
$ws.Range("K55").Value = 5149.5
$ws.Range("L55").Value = 356.4
$ws.Range("M55").Value = -4935.5
$ws.Range("N55").Value = -784.4
$ws.Range("H70").Value = 1225.4375
$ws.Range("I70").Value = 687.6667
$ws.Range("J70").Value = 1548.1
$ws.Range("K70").Value = 2063.0001
$ws.Range("L70").Value = 4644.299999999999
$ws.Range("M70").Value = -1793.0001
$ws.Range("N70").Value = -5184.299999999999
$ws.Range("H73").Value = 1225.4375
$ws.Range("I73").Value = 687.6667
$ws.Range("J73").Value = 1548.1
$ws.Range("K73").Value = 2063.0001
$ws.Range("L73").Value = 4644.299999999999
$ws.Range("M73").Value = -1127.0001
$ws.Range("N73").Value = -6516.299999999999
$ws.Range("H92").Value = 2757.5715
$ws.Range("I92").Value = 2755.45
$ws.Range("K92").Value = 2755.45
$ws.Range("M92").Value = -1507.45
$ws.Range("H98").Value = 1512.8125
$ws.Range("I98").Value = 1554.2307
$ws.Range("J98").Value = 1333.3334
$ws.Range("K98").Value = 1554.2307
$ws.Range("L98").Value = 1333.3334
$ws.Range("M98").Value = -56.23070000000007
$ws.Range("N98").Value = -4329.3334
$ws.Range("H122").Value = 1512.8125
$ws.Range("I122").Value = 1554.2307
$ws.Range("J122").Value = 1333.3334
$ws.Range("K122").Value = 4662.6921
$ws.Range("L122").Value = 4000.0002
$ws.Range("M122").Value = -2212.6921
$ws.Range("N122").Value = -8900.0002
$ws.Range("H133").Value = 47680
$ws.Range("J133").Value = 47680
$ws.Range("L133").Value = 47680
$ws.Range("N133").Value = -57800
$ws.Range("H137").Value = 1149.8438
$ws.Range("I137").Value = 971.26086
$ws.Range("J137").Value = 1606.2222
$ws.Range("K137").Value = 2913.78258
$ws.Range("L137").Value = 4818.6666
$ws.Range("M137").Value = -363.7825800000001
$ws.Range("N137").Value = -9918.6666
$ws.Range("H138").Value = 4247.83
$ws.Range("I138").Value = 889.8889
$ws.Range("J138").Value = 4984.939
$ws.Range("K138").Value = 2669.6667
$ws.Range("L138").Value = 14954.817
$ws.Range("M138").Value = 2470.3333
$ws.Range("N138").Value = -25234.817

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1187.4667
$ws.Range("I61").Value = 812.44446
$ws.Range("J61").Value = 1750
$ws.Range("K61").Value = 812.44446
$ws.Range("L61").Value = 1750
$ws.Range("M61").Value = -600.44446
$ws.Range("N61").Value = -2174
$ws.Range("H102").Value = 2742.9333
$ws.Range("I102").Value = 2707.3076
$ws.Range("K102").Value = 2707.3076
$ws.Range("M102").Value = -1085.3076
$ws.Range("H132").Value = 1541.9778
$ws.Range("I132").Value = 1177.5834
$ws.Range("J132").Value = 2999.5557
$ws.Range("K132").Value = 3532.7502
$ws.Range("L132").Value = 8998.667099999999
$ws.Range("M132").Value = -1002.7502
$ws.Range("N132").Value = -14058.6671
$ws.Range("H136").Value = 1187.4667
$ws.Range("I136").Value = 812.44446
$ws.Range("J136").Value = 1750
$ws.Range("K136").Value = 2437.33338
$ws.Range("L136").Value = 5250
$ws.Range("M136").Value = 112.66662
$ws.Range("N136").Value = -10350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 824.25
$ws.Range("I64").Value = 860
$ws.Range("J64").Value = 808
$ws.Range("K64").Value = 860
$ws.Range("L64").Value = 808
$ws.Range("M64").Value = -635
$ws.Range("N64").Value = -1258
$ws.Range("H67").Value = 824.25
$ws.Range("I67").Value = 860
$ws.Range("J67").Value = 808
$ws.Range("K67").Value = 860
$ws.Range("L67").Value = 808
$ws.Range("M67").Value = -80
$ws.Range("N67").Value = -2368
$ws.Range("H105").Value = 10353.087
$ws.Range("I105").Value = 9247.647000000001
$ws.Range("J105").Value = 13485.167
$ws.Range("K105").Value = 9247.647000000001
$ws.Range("L105").Value = 13485.167
$ws.Range("M105").Value = -7500.647000000001
$ws.Range("N105").Value = -16979.167
$ws.Range("H134").Value = 20596.906
$ws.Range("I134").Value = 1676.9736
$ws.Range("J134").Value = 68527.39999999999
$ws.Range("K134").Value = 5030.9208
$ws.Range("L134").Value = 205582.2
$ws.Range("M134").Value = -2495.9208
$ws.Range("N134").Value = -210652.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2807.1667
$ws.Range("I31").Value = 2789.5527
$ws.Range("J31").Value = 2974.5
$ws.Range("K31").Value = 2789.5527
$ws.Range("L31").Value = 2974.5
$ws.Range("M31").Value = -2494.5527
$ws.Range("N31").Value = -3564.5
$ws.Range("H34").Value = 2807.1667
$ws.Range("I34").Value = 2789.5527
$ws.Range("J34").Value = 2974.5
$ws.Range("K34").Value = 2789.5527
$ws.Range("L34").Value = 2974.5
$ws.Range("M34").Value = -2587.5527
$ws.Range("N34").Value = -3378.5
$ws.Range("H134").Value = 1412.6818
$ws.Range("I134").Value = 483
$ws.Range("J134").Value = 2755.5557
$ws.Range("K134").Value = 1449
$ws.Range("L134").Value = 8266.667099999999
$ws.Range("M134").Value = 1086
$ws.Range("N134").Value = -13336.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 368.53845
$ws.Range("J38").Value = 622.8570999999999
$ws.Range("L38").Value = 1868.5713
$ws.Range("N38").Value = -2562.5713
$ws.Range("H139").Value = 1657.7778
$ws.Range("I139").Value = 1721.25
$ws.Range("J139").Value = 1150
$ws.Range("K139").Value = 5163.75
$ws.Range("L139").Value = 3450
$ws.Range("M139").Value = -23.75
$ws.Range("N139").Value = -13730

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3738.0952
$ws.Range("I80").Value = 3676.4707
$ws.Range("K80").Value = 3676.4707
$ws.Range("M80").Value = -2678.4707
$ws.Range("H83").Value = 3738.0952
$ws.Range("I83").Value = 3676.4707
$ws.Range("K83").Value = 18382.3535
$ws.Range("M83").Value = -13390.3535
$ws.Range("H97").Value = 2938.84
$ws.Range("J97").Value = 1997.3636
$ws.Range("L97").Value = 1997.3636
$ws.Range("N97").Value = -2989.3636
$ws.Range("H102").Value = 2369.05
$ws.Range("I102").Value = 1691.5
$ws.Range("K102").Value = 1691.5
$ws.Range("M102").Value = -69.5
$ws.Range("H122").Value = 1317274.1
$ws.Range("I122").Value = 1880977.2
$ws.Range("J122").Value = 1966.6666
$ws.Range("K122").Value = 5642931.6
$ws.Range("L122").Value = 5899.9998
$ws.Range("M122").Value = -5640481.6
$ws.Range("N122").Value = -10799.9998
$ws.Range("H126").Value = 2343.0435
$ws.Range("I126").Value = 2240.1428
$ws.Range("J126").Value = 2503.111
$ws.Range("K126").Value = 6720.428400000001
$ws.Range("L126").Value = 7509.333
$ws.Range("M126").Value = -4250.428400000001
$ws.Range("N126").Value = -12449.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18521400
$ws.Range("I7").Value = 3458
$ws.Range("J7").Value = 111111110
$ws.Range("K7").Value = 3458
$ws.Range("L7").Value = 111111110
$ws.Range("M7").Value = -3346
$ws.Range("N7").Value = -111111334
$ws.Range("H40").Value = 919732.75
$ws.Range("I40").Value = 1684835
$ws.Range("K40").Value = 1684835
$ws.Range("M40").Value = -1684699
$ws.Range("H122").Value = 16999.166
$ws.Range("I122").Value = 16999.166
$ws.Range("K122").Value = 50997.49800000001
$ws.Range("M122").Value = -48547.49800000001
$ws.Range("H126").Value = 18521400
$ws.Range("I126").Value = 3458
$ws.Range("J126").Value = 111111110
$ws.Range("K126").Value = 10374
$ws.Range("L126").Value = 333333330
$ws.Range("M126").Value = -7904
$ws.Range("N126").Value = -333338270
$ws.Range("H136").Value = 1737.8485
$ws.Range("I136").Value = 1191.5094
$ws.Range("J136").Value = 3965.2307
$ws.Range("K136").Value = 3574.5282
$ws.Range("L136").Value = 11895.6921
$ws.Range("M136").Value = -1024.5282
$ws.Range("N136").Value = -16995.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 863.75
$ws.Range("I122").Value = 825
$ws.Range("K122").Value = 2475
$ws.Range("M122").Value = -25
$ws.Range("H132").Value = 1247.9783
$ws.Range("I132").Value = 964.54285
$ws.Range("J132").Value = 2149.818
$ws.Range("K132").Value = 2893.62855
$ws.Range("L132").Value = 6449.454000000001
$ws.Range("M132").Value = -363.6285500000004
$ws.Range("N132").Value = -11509.454
